$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (sheet name in workbook.xml: "test1" -> "validLogin")
$ws.Name = "validLogin"

# Set up the login credential headers/values
# Order of assignment controls the shared-string table order:
# username(0), manager(1), password(2), admin(3)
$ws.Range("A1").Value = "username"
$ws.Range("B2").Value = "manager"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"

# Select A2 as the active cell, matching the saved selection state
$ws.Range("A2").Select()
